$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 3562.5
$ws.Range("I70").Value = 3357.1428
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 10071.4284
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -9801.428400000001
$ws.Range("N70").Value = -15540

# Row 73
$ws.Range("H73").Value = 3562.5
$ws.Range("I73").Value = 3357.1428
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 10071.4284
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -9135.428400000001
$ws.Range("N73").Value = -16872

# Row 103
$ws.Range("H103").Value = 1376
$ws.Range("I103").Value = 1002
$ws.Range("K103").Value = 3006
$ws.Range("M103").Value = -2420

# Row 135
$ws.Range("H135").Value = 483
$ws.Range("I135").Value = 483
$ws.Range("K135").Value = 4347
$ws.Range("M135").Value = -1812

# Row 137
$ws.Range("H137").Value = 966
$ws.Range("I137").Value = 966
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2898
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -348
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 13
$ws.Range("H13").Value = 9166.333000000001
$ws.Range("I13").Value = 7500
$ws.Range("K13").Value = 7500
$ws.Range("M13").Value = -7356

# Row 14
$ws.Range("H14").Value = 200630
$ws.Range("J14").Value = 250
$ws.Range("L14").Value = 250
$ws.Range("N14").Value = -600

# Row 32
$ws.Range("H32").Value = 2891.15
$ws.Range("I32").Value = 2891.15
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2891.15
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2604.15
$ws.Range("N32").ClearContents()

# Row 61
$ws.Range("H61").Value = 1950
$ws.Range("I61").Value = 1900
$ws.Range("K61").Value = 1900
$ws.Range("M61").Value = -1688

# Row 76
$ws.Range("H76").Value = 32483.75
$ws.Range("J76").Value = 32483.75
$ws.Range("L76").Value = 32483.75
$ws.Range("N76").Value = -33159.75

# Row 79
$ws.Range("H79").Value = 32483.75
$ws.Range("J79").Value = 32483.75
$ws.Range("L79").Value = 32483.75
$ws.Range("N79").Value = -34823.75

# Row 122
$ws.Range("H122").Value = 4814.4614
$ws.Range("I122").Value = 5080.636
$ws.Range("J122").Value = 3350.5
$ws.Range("K122").Value = 15241.908
$ws.Range("L122").Value = 10051.5
$ws.Range("M122").Value = -12791.908
$ws.Range("N122").Value = -14951.5

# Row 136
$ws.Range("H136").Value = 1950
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -3150

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 54
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

# Row 12
$ws.Range("H12").Value = 1131.8889
$ws.Range("I12").Value = 633.8570999999999
$ws.Range("K12").Value = 633.8570999999999
$ws.Range("M12").Value = -465.8570999999999

# Row 16
$ws.Range("H16").Value = 1595.8
$ws.Range("J16").Value = 2393
$ws.Range("L16").Value = 2393
$ws.Range("N16").Value = -2733

$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 3640.3076
$ws.Range("I8").Value = 1076.75
$ws.Range("K8").Value = 1076.75
$ws.Range("M8").Value = -936.75

# Row 14
$ws.Range("H14").Value = 10
$ws.Range("I14").Value = 10
$ws.Range("K14").Value = 10
$ws.Range("M14").Value = 160

# Row 33
$ws.Range("H33").Value = 25979.8
$ws.Range("I33").Value = 3300
$ws.Range("K33").Value = 3300
$ws.Range("M33").Value = -2921

# Row 36
$ws.Range("H36").Value = 32498
$ws.Range("I36").Value = 19998.8
$ws.Range("J36").Value = 53330
$ws.Range("K36").Value = 19998.8
$ws.Range("L36").Value = 53330
$ws.Range("M36").Value = -19610.8
$ws.Range("N36").Value = -54106

# Row 40
$ws.Range("H40").Value = 32498
$ws.Range("I40").Value = 19998.8
$ws.Range("J40").Value = 53330
$ws.Range("K40").Value = 19998.8
$ws.Range("L40").Value = 53330
$ws.Range("M40").Value = -19838.8
$ws.Range("N40").Value = -53650

# Row 55
$ws.Range("H55").Value = 49999.25
$ws.Range("I55").Value = 59999
$ws.Range("J55").Value = 39999.5
$ws.Range("K55").Value = 59999
$ws.Range("L55").Value = 39999.5
$ws.Range("M55").Value = -59684
$ws.Range("N55").Value = -40629.5

# Row 58
$ws.Range("H58").Value = 20500
$ws.Range("I58").Value = 20500
$ws.Range("K58").Value = 20500
$ws.Range("M58").Value = -20297

# Row 86
$ws.Range("H86").Value = 1000000000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 1000000000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 134
$ws.Range("H134").Value = 2708.7
$ws.Range("I134").Value = 2708.7
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8126.099999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5591.099999999999
$ws.Range("N134").ClearContents()

# Row 136
$ws.Range("H136").Value = 20500
$ws.Range("I136").Value = 20500
$ws.Range("K136").Value = 61500
$ws.Range("M136").Value = -58950

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 13
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("K13").Value = 1500
$ws.Range("M13").Value = -1332

# Row 29
$ws.Range("H29").Value = 159.11111
$ws.Range("I29").Value = 280.5
$ws.Range("J29").Value = 62
$ws.Range("K29").Value = 841.5
$ws.Range("L29").Value = 186
$ws.Range("M29").Value = -564.5
$ws.Range("N29").Value = -740

# Row 46
$ws.Range("H46").Value = 2189.8823
$ws.Range("I46").Value = 629
$ws.Range("K46").Value = 1887
$ws.Range("M46").Value = -1796

# Row 134
$ws.Range("H134").Value = 987.75
$ws.Range("I134").Value = 987.75
$ws.Range("K134").Value = 2963.25
$ws.Range("M134").Value = 2106.75

# Row 136
$ws.Range("H136").Value = 1030
$ws.Range("I136").Value = 1030
$ws.Range("K136").Value = 3090
$ws.Range("M136").Value = 2010

# Row 139
$ws.Range("H139").Value = 4694.8
$ws.Range("I139").Value = 4743.5
$ws.Range("K139").Value = 14230.5
$ws.Range("M139").Value = -9090.5

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 6166.6665
$ws.Range("I3").Value = 3500
$ws.Range("K3").Value = 3500
$ws.Range("M3").Value = -3384

# Row 14
$ws.Range("H14").Value = 16667367
$ws.Range("J14").Value = 1050
$ws.Range("L14").Value = 1050
$ws.Range("N14").Value = -1386

# Row 124
$ws.Range("H124").Value = 99999.75
$ws.Range("J124").Value = 99999.75
$ws.Range("L124").Value = 99999.75
$ws.Range("N124").Value = -109819.75

$ws = $wb.Worksheets.Item("LTW")
# Row 41
$ws.Range("H41").Value = 35000
$ws.Range("J41").Value = 35000
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35876

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 6838.75
$ws.Range("J6").Value = 3451.6667
$ws.Range("L6").Value = 3451.6667
$ws.Range("N6").Value = -3681.6667

# Row 9
$ws.Range("H9").Value = 999
$ws.Range("I9").Value = 999
$ws.Range("K9").Value = 999
$ws.Range("M9").Value = -859

# Row 13
$ws.Range("H13").Value = 1005
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 81
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 3500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5939
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 3500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 35000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -29696
$ws.Range("N84").ClearContents()

